$wb = $excel.ActiveWorkbook

$updates1 = @{
    2  = 919
    7  = 4520
    8  = 2675
    10 = 2639
    14 = 1684
    15 = 698
    16 = 360
    26 = 595
    30 = 461
    32 = 1270
    35 = 1316
    36 = 2175
    37 = 332
    39 = 571
    41 = 40
    43 = 718
    44 = 1399
    45 = 155
    47 = 458
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

$updates4 = @{
    2  = 919
    5  = 4520
    6  = 2675
    7  = 2639
    8  = 1684
    11 = 698
    12 = 360
    21 = 595
    27 = 461
    29 = 1270
    34 = 2175
    35 = 332
    39 = 571
    41 = 40
    43 = 718
    44 = 1399
    46 = 155
    47 = 458
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
